$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (009540.KS / HDKSOE)
$ws.Range("D2").Value = 426000
$ws.Range("E2").Value = 46.2
$ws.Range("F2").Value = 3.9
$ws.Range("N2").Value = 54.84087454262382

# Row 3 (010620.KS / HD HYUNDAI MIPO)
$ws.Range("N3").Value = 54.84087454262382

# Row 4 (042660.KS / Hanwha Ocean)
$ws.Range("D4").Value = 106700
$ws.Range("E4").Value = 17.1
$ws.Range("F4").Value = -1.02
$ws.Range("N4").Value = 54.84087454262382

# Row 5 (010140.KS / SamsungHvyInd)
$ws.Range("D5").Value = 24850
$ws.Range("E5").Value = 39
$ws.Range("F5").Value = 1.02
$ws.Range("I5").Value = 63
$ws.Range("K5").Value = 47.7
$ws.Range("N5").Value = 54.84087454262382
